$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.987640333333333
$ws.Range("H2").Value = 5.962921
$ws.Range("I2").Value = 0.0278174819837782
$ws.Range("J2").Value = 0.0278174819837782
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.071057
$ws.Range("N2").Value = 102.213171
$ws.Range("O2").Value = 0.5537562116045693
$ws.Range("P2").Value = 0.5537562116045693
$ws.Range("Q2").Value = 67.72100709249901
$ws.Range("R2").Value = 609.489063832491
$ws.Range("S2").Value = 0.01540410343971538
$ws.Range("T2").Value = 0.01540410343971538

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.987640333333333
$ws.Range("H3").Value = 5.962921
$ws.Range("I3").Value = 0.0278174819837782
$ws.Range("J3").Value = 0.0278174819837782
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.28977566666667
$ws.Range("N3").Value = 57.869327
$ws.Range("O3").Value = 0.3135163401556734
$ws.Range("P3").Value = 0.3135163401556735
$ws.Range("Q3").Value = 38.34113613601856
$ws.Range("R3").Value = 345.070225224167
$ws.Range("S3").Value = 0.008721235143900524
$ws.Range("T3").Value = 0.008721235143900526

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.987640333333333
$ws.Range("H4").Value = 5.962921
$ws.Range("I4").Value = 0.0278174819837782
$ws.Range("J4").Value = 0.0278174819837782
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.368545333333333
$ws.Range("N4").Value = 19.105636
$ws.Range("O4").Value = 0.1035078406055505
$ws.Range("P4").Value = 0.1035078406055505
$ws.Range("Q4").Value = 12.65837756919511
$ws.Range("R4").Value = 113.925398122756
$ws.Range("S4").Value = 0.002879327491224686
$ws.Range("T4").Value = 0.002879327491224687

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.987640333333333
$ws.Range("H5").Value = 5.962921
$ws.Range("I5").Value = 0.0278174819837782
$ws.Range("J5").Value = 0.0278174819837782
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.7978
$ws.Range("N5").Value = 5.3934
$ws.Range("O5").Value = 0.02921960763420679
$ws.Range("P5").Value = 0.02921960763420679
$ws.Range("Q5").Value = 3.573379791266666
$ws.Range("R5").Value = 32.1604181214
$ws.Range("S5").Value = 0.0008128159089376153
$ws.Range("T5").Value = 0.0008128159089376155

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 56.81334033333334
$ws.Range("H6").Value = 170.440021
$ws.Range("I6").Value = 0.7951157181995667
$ws.Range("J6").Value = 0.7951157181995667
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 34.071057
$ws.Range("N6").Value = 102.213171
$ws.Range("O6").Value = 0.5537562116045693
$ws.Range("P6").Value = 0.5537562116045693
$ws.Range("Q6").Value = 1935.690556857399
$ws.Range("R6").Value = 17421.21501171659
$ws.Range("S6").Value = 0.4403002678974384
$ws.Range("T6").Value = 0.4403002678974384

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 56.81334033333334
$ws.Range("H7").Value = 170.440021
$ws.Range("I7").Value = 0.7951157181995667
$ws.Range("J7").Value = 0.7951157181995667
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.28977566666667
$ws.Range("N7").Value = 57.869327
$ws.Range("O7").Value = 0.3135163401556734
$ws.Range("P7").Value = 0.3135163401556735
$ws.Range("Q7").Value = 1095.916589903985
$ws.Range("R7").Value = 9863.249309135866
$ws.Range("S7").Value = 0.2492817699701779
$ws.Range("T7").Value = 0.249281769970178

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 56.81334033333334
$ws.Range("H8").Value = 170.440021
$ws.Range("I8").Value = 0.7951157181995667
$ws.Range("J8").Value = 0.7951157181995667
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.368545333333333
$ws.Range("N8").Value = 19.105636
$ws.Range("O8").Value = 0.1035078406055505
$ws.Range("P8").Value = 0.1035078406055505
$ws.Range("Q8").Value = 361.8183334509284
$ws.Range("R8").Value = 3256.365001058356
$ws.Range("S8").Value = 0.08230071102236855
$ws.Range("T8").Value = 0.08230071102236856

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 56.81334033333334
$ws.Range("H9").Value = 170.440021
$ws.Range("I9").Value = 0.7951157181995667
$ws.Range("J9").Value = 0.7951157181995667
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.7978
$ws.Range("N9").Value = 5.3934
$ws.Range("O9").Value = 0.02921960763420679
$ws.Range("P9").Value = 0.02921960763420679
$ws.Range("Q9").Value = 102.1390232512667
$ws.Range("R9").Value = 919.2512092613999
$ws.Range("S9").Value = 0.02323296930958187
$ws.Range("T9").Value = 0.02323296930958188

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7501196666666666
$ws.Range("H10").Value = 2.250359
$ws.Range("I10").Value = 0.01049809664416703
$ws.Range("J10").Value = 0.01049809664416703
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 34.071057
$ws.Range("N10").Value = 102.213171
$ws.Range("O10").Value = 0.5537562116045693
$ws.Range("P10").Value = 0.5537562116045693
$ws.Range("Q10").Value = 25.557369919821
$ws.Range("R10").Value = 230.016329278389
$ws.Range("S10").Value = 0.005813386226732578
$ws.Range("T10").Value = 0.005813386226732579

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.7501196666666666
$ws.Range("H11").Value = 2.250359
$ws.Range("I11").Value = 0.01049809664416703
$ws.Range("J11").Value = 0.01049809664416703
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 19.28977566666667
$ws.Range("N11").Value = 57.869327
$ws.Range("O11").Value = 0.3135163401556734
$ws.Range("P11").Value = 0.3135163401556735
$ws.Range("Q11").Value = 14.46964009315478
$ws.Range("R11").Value = 130.226760838393
$ws.Range("S11").Value = 0.003291324838479805
$ws.Range("T11").Value = 0.003291324838479806

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.7501196666666666
$ws.Range("H12").Value = 2.250359
$ws.Range("I12").Value = 0.01049809664416703
$ws.Range("J12").Value = 0.01049809664416703
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.368545333333333
$ws.Range("N12").Value = 19.105636
$ws.Range("O12").Value = 0.1035078406055505
$ws.Range("P12").Value = 0.1035078406055505
$ws.Range("Q12").Value = 4.777171102591555
$ws.Range("R12").Value = 42.99453992332399
$ws.Range("S12").Value = 0.001086635314106106
$ws.Range("T12").Value = 0.001086635314106106

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.7501196666666666
$ws.Range("H13").Value = 2.250359
$ws.Range("I13").Value = 0.01049809664416703
$ws.Range("J13").Value = 0.01049809664416703
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.7978
$ws.Range("N13").Value = 5.3934
$ws.Range("O13").Value = 0.02921960763420679
$ws.Range("P13").Value = 0.02921960763420679
$ws.Range("Q13").Value = 1.348565136733333
$ws.Range("R13").Value = 12.1370862306
$ws.Range("S13").Value = 0.0003067502648485437
$ws.Range("T13").Value = 0.0003067502648485438

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 11.90182033333333
$ws.Range("H14").Value = 35.705461
$ws.Range("I14").Value = 0.166568703172488
$ws.Range("J14").Value = 0.166568703172488
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 34.071057
$ws.Range("N14").Value = 102.213171
$ws.Range("O14").Value = 0.5537562116045693
$ws.Range("P14").Value = 0.5537562116045693
$ws.Range("Q14").Value = 405.507598980759
$ws.Range("R14").Value = 3649.568390826831
$ws.Range("S14").Value = 0.09223845404068295
$ws.Range("T14").Value = 0.09223845404068295

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 11.90182033333333
$ws.Range("H15").Value = 35.705461
$ws.Range("I15").Value = 0.166568703172488
$ws.Range("J15").Value = 0.166568703172488
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 19.28977566666667
$ws.Range("N15").Value = 57.869327
$ws.Range("O15").Value = 0.3135163401556734
$ws.Range("P15").Value = 0.3135163401556735
$ws.Range("Q15").Value = 229.5834442549719
$ws.Range("R15").Value = 2066.250998294747
$ws.Range("S15").Value = 0.05222201020311514
$ws.Range("T15").Value = 0.05222201020311515

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 11.90182033333333
$ws.Range("H16").Value = 35.705461
$ws.Range("I16").Value = 0.166568703172488
$ws.Range("J16").Value = 0.166568703172488
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.368545333333333
$ws.Range("N16").Value = 19.105636
$ws.Range("O16").Value = 0.1035078406055505
$ws.Range("P16").Value = 0.1035078406055505
$ws.Range("Q16").Value = 75.79728234202177
$ws.Range("R16").Value = 682.1755410781959
$ws.Range("S16").Value = 0.01724116677785114
$ws.Range("T16").Value = 0.01724116677785114

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 11.90182033333333
$ws.Range("H17").Value = 35.705461
$ws.Range("I17").Value = 0.166568703172488
$ws.Range("J17").Value = 0.166568703172488
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.7978
$ws.Range("N17").Value = 5.3934
$ws.Range("O17").Value = 0.02921960763420679
$ws.Range("P17").Value = 0.02921960763420679
$ws.Range("Q17").Value = 21.39709259526666
$ws.Range("R17").Value = 192.5738333574
$ws.Range("S17").Value = 0.004867072150838754
$ws.Range("T17").Value = 0.004867072150838755
